# Apply updated crypto price/volume values to Sheet1 (cells D2:E51).
# Values are stored as text in the original workbook (t="inlineStr"), so for
# any cell whose new value would otherwise be auto-detected by Excel as a
# number, we force the cell to Text format first, then assign the string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.497.09"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.954.01"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.03"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.21"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0788"
$ws.Range("E10").Value = "  -7.94%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.04"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").Value = "2.241.48"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "1.964.68"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "36.424.26"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.71"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.19"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.04"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.13"
$ws.Range("E27").Value = "  +6.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.19"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.36"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.48"
$ws.Range("E35").Value = "  +12.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("E39").Value = "  -14.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0978"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0210"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "1.368.77"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.77"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.83"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.13"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "2.131.71"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.91"
$ws.Range("E51").Value = "  -1.19%  "
